$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 13:35"

# Row 14 - India
$ws.Range("B14").Value = 102231
$ws.Range("C14").Value = 1903
$ws.Range("D14").Value = 39658
$ws.Range("E14").Value = 59404
$ws.Range("G14").Value = 13
$ws.Range("H14").Value = 3169

# Row 75 - Uzbekistan
$ws.Range("B75").Value = 2825
$ws.Range("C75").Value = 34
$ws.Range("D75").Value = 2338
$ws.Range("E75").Value = 474

# Row 80 - Bosnia y Herzegovina
$ws.Range("B80").Value = 2321
$ws.Range("C80").Value = 17
$ws.Range("D80").Value = 1522
$ws.Range("E80").Value = 665
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 134

# Row 106 - Sri Lanka
$ws.Range("D106").Value = 569
$ws.Range("E106").Value = 414
